$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unnecessary "nome" (name) column B entirely; all columns to
# the right shift left by one.
$ws.Columns("B").Delete()

# Put selection on B1 to match the saved view state.
$ws.Range("B1").Select()
